$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4.259755363148439
$ws.Range("D2").Value = 6.423941403031091
$ws.Range("E2").Value = 9.677591315237592
$ws.Range("F2").Value = 35.63042287516749
$ws.Range("G2").Value = 3.669130531475786
$ws.Range("I2").Value = 28.52160734120368
$ws.Range("K2").Value = 15.42243458850324
$ws.Range("M2").Value = 16.48296303869194

# Row 3
$ws.Range("C3").Value = 4.262655045587705
$ws.Range("D3").Value = 6.445699409744611
$ws.Range("E3").Value = 9.559697735725409
$ws.Range("F3").Value = 35.27445817448563
$ws.Range("G3").Value = 3.673175339344439
$ws.Range("I3").Value = 28.31762919710955
$ws.Range("K3").Value = 15.0203023777068
$ws.Range("M3").Value = 16.23935062608864

# Row 4
$ws.Range("C4").Value = 4.264395839487481
$ws.Range("D4").Value = 6.45983434745262
$ws.Range("E4").Value = 9.489821267916994
$ws.Range("F4").Value = 35.06544192608423
$ws.Range("G4").Value = 3.675784285579228
$ws.Range("I4").Value = 28.20013690018182
$ws.Range("K4").Value = 14.7727986689093
$ws.Range("M4").Value = 16.09310557048652

# Row 5
$ws.Range("C5").Value = 4.265095631974321
$ws.Range("D5").Value = 6.465788709175206
$ws.Range("E5").Value = 9.462006017182526
$ws.Range("F5").Value = 34.9827416512185
$ws.Range("G5").Value = 3.676879119788919
$ws.Range("I5").Value = 28.15423807520764
$ws.Range("K5").Value = 14.67196163302021
$ws.Range("M5").Value = 16.0344246150691

# Row 6
$ws.Range("C6").Value = 4.265211263302837
$ws.Range("D6").Value = 6.466789138241885
$ws.Range("E6").Value = 9.45742798006791
$ws.Range("F6").Value = 34.96916091352401
$ws.Range("G6").Value = 3.677062832625325
$ws.Range("I6").Value = 28.14673704202137
$ws.Range("K6").Value = 14.65522424287422
$ws.Range("M6").Value = 16.02473816343275

# Row 7
$ws.Range("C7").Value = 4.264405315522703
$ws.Range("D7").Value = 6.459913864496175
$ws.Range("E7").Value = 9.489443432903936
$ws.Range("F7").Value = 35.06431648395673
$ws.Range("G7").Value = 3.675798922505728
$ws.Range("I7").Value = 28.19950983674332
$ws.Range("K7").Value = 14.77143842148365
$ws.Range("M7").Value = 16.09231037341072

# Row 8
$ws.Range("C8").Value = 4.260763618244028
$ws.Range("D8").Value = 6.431282180304247
$ws.Range("E8").Value = 9.636438219132939
$ws.Range("F8").Value = 35.50575133013155
$ws.Range("G8").Value = 3.670499229287074
$ws.Range("I8").Value = 28.44968685309243
$ws.Range("K8").Value = 15.28400338970373
$ws.Range("M8").Value = 16.39832005885795

# Row 9
$ws.Range("C9").Value = 4.253291291454218
$ws.Range("D9").Value = 6.381321046820116
$ws.Range("E9").Value = 9.943280420754993
$ws.Range("F9").Value = 36.4435962405865
$ws.Range("G9").Value = 3.661095736715591
$ws.Range("I9").Value = 29.00022844744676
$ws.Range("K9").Value = 16.27726799447985
$ws.Range("M9").Value = 17.02141872036263

# Row 10
$ws.Range("C10").Value = 4.247576683826805
$ws.Range("D10").Value = 6.348430816790744
$ws.Range("E10").Value = 10.17817545829172
$ws.Range("F10").Value = 37.17153053207405
$ws.Range("G10").Value = 3.654781691825137
$ws.Range("I10").Value = 29.43892398043203
$ws.Range("K10").Value = 16.99080538426381
$ws.Range("M10").Value = 17.48859783379131

# Row 11
$ws.Range("C11").Value = 4.244923303016866
$ws.Range("D11").Value = 6.33430710320345
$ws.Range("E11").Value = 10.28667475212068
$ws.Range("F11").Value = 37.50996707320267
$ws.Range("G11").Value = 3.652036616058714
$ws.Range("I11").Value = 29.64536344440646
$ws.Range("K11").Value = 17.31019119399416
$ws.Range("M11").Value = 17.70219185564341

# Row 12
$ws.Range("C12").Value = 4.243910421185096
$ws.Range("D12").Value = 6.329080283426671
$ws.Range("E12").Value = 10.32796188240006
$ws.Range("F12").Value = 37.63907162305646
$ws.Range("G12").Value = 3.651015282557911
$ws.Range("I12").Value = 29.72447377668333
$ws.Range("K12").Value = 17.43025846730401
$ws.Range("M12").Value = 17.78314885812443

# Row 13
$ws.Range("C13").Value = 4.244128929716906
$ws.Range("D13").Value = 6.330200551660788
$ws.Range("E13").Value = 10.31906159600278
$ws.Range("F13").Value = 37.61122623211931
$ws.Range("G13").Value = 3.651234439021937
$ws.Range("I13").Value = 29.70739517072478
$ws.Range("K13").Value = 17.4044407781292
$ws.Range("M13").Value = 17.7657113053251

# Row 14
$ws.Range("C14").Value = 4.244840137156096
$ws.Range("D14").Value = 6.333874648292537
$ws.Range("E14").Value = 10.29006764866646
$ws.Range("F14").Value = 37.52057021907081
$ws.Range("G14").Value = 3.651952226963548
$ws.Range("I14").Value = 29.65185344416076
$ws.Range("K14").Value = 17.32008747940687
$ws.Range("M14").Value = 17.70885121576126

# Row 15
$ws.Range("C15").Value = 4.245274705880371
$ws.Range("D15").Value = 6.336140994492926
$ws.Range("E15").Value = 10.27233311245553
$ws.Range("F15").Value = 37.46516090259831
$ws.Range("G15").Value = 3.652394254954119
$ws.Range("I15").Value = 29.61795286046134
$ws.Range("K15").Value = 17.26830077881539
$ws.Range("M15").Value = 17.67402999853784

# Row 16
$ws.Range("C16").Value = 4.247748974595984
$ws.Range("D16").Value = 6.349370785605593
$ws.Range("E16").Value = 10.17111481409788
$ws.Range("F16").Value = 37.14955093443798
$ws.Range("G16").Value = 3.654963638076272
$ws.Range("I16").Value = 29.42556658507429
$ws.Range("K16").Value = 16.96981677424199
$ws.Range("M16").Value = 17.47465359655493

# Row 17
$ws.Range("C17").Value = 4.24925281989448
$ws.Range("D17").Value = 6.357702215279022
$ws.Range("E17").Value = 10.1094164544044
$ws.Range("F17").Value = 36.9577301688736
$ws.Range("G17").Value = 3.656572364308778
$ws.Range("I17").Value = 29.30926805764152
$ws.Range("K17").Value = 16.78527821222676
$ws.Range("M17").Value = 17.35255919568446

# Row 18
$ws.Range("C18").Value = 4.250112752192464
$ws.Range("D18").Value = 6.362573069128034
$ws.Range("E18").Value = 10.07408610589109
$ws.Range("F18").Value = 36.84809293198159
$ws.Range("G18").Value = 3.657509642696023
$ws.Range("I18").Value = 29.24302658792333
$ws.Range("K18").Value = 16.67865350756567
$ws.Range("M18").Value = 17.28243932011781

# Row 19
$ws.Range("C19").Value = 4.250403056246687
$ws.Range("D19").Value = 6.364235767035445
$ws.Range("E19").Value = 10.06215189433169
$ws.Range("F19").Value = 36.81109378538096
$ws.Range("G19").Value = 3.657829050754195
$ws.Range("I19").Value = 29.22071162460319
$ws.Range("K19").Value = 16.64247335367502
$ws.Range("M19").Value = 17.25871856005887

# Row 20
$ws.Range("C20").Value = 4.249093257285457
$ws.Range("D20").Value = 6.356807154093941
$ws.Range("E20").Value = 10.11596836837092
$ws.Range("F20").Value = 36.97807883677137
$ws.Range("G20").Value = 3.656399873519119
$ws.Range("I20").Value = 29.32158129630449
$ws.Range("K20").Value = 16.80497361221786
$ws.Range("M20").Value = 17.36554601485985

# Row 21
$ws.Range("C21").Value = 4.24463146117792
$ws.Range("D21").Value = 6.332792170874793
$ws.Range("E21").Value = 10.29857870776538
$ws.Range("F21").Value = 37.54717324282657
$ws.Range("G21").Value = 3.651740903200421
$ws.Range("I21").Value = 29.66814241586358
$ws.Range("K21").Value = 17.34488886623199
$ws.Range("M21").Value = 17.72555103632238

# Row 22
$ws.Range("C22").Value = 4.241668011985817
$ws.Range("D22").Value = 6.317805808442952
$ws.Range("E22").Value = 10.41907914147024
$ws.Range("F22").Value = 37.9245745665591
$ws.Range("G22").Value = 3.64880183016684
$ws.Range("I22").Value = 29.90007122904572
$ws.Range("K22").Value = 17.69258964720095
$ws.Range("M22").Value = 17.96122832742732

# Row 23
$ws.Range("C23").Value = 4.243254120354056
$ws.Range("D23").Value = 6.325739119900284
$ws.Range("E23").Value = 10.35467193976282
$ws.Range("F23").Value = 37.72268269791065
$ws.Range("G23").Value = 3.650360826654901
$ws.Range("I23").Value = 29.77580735351924
$ws.Range("K23").Value = 17.50752749584014
$ws.Range("M23").Value = 17.83543295925401

# Row 24
$ws.Range("C24").Value = 4.249165409991713
$ws.Range("D24").Value = 6.357211559020906
$ws.Range("E24").Value = 10.11300580661239
$ws.Range("F24").Value = 36.96887719317209
$ws.Range("G24").Value = 3.656477817971461
$ws.Range("I24").Value = 29.31601254567671
$ws.Range("K24").Value = 16.79607096603843
$ws.Range("M24").Value = 17.35967443783675

# Row 25
$ws.Range("C25").Value = 4.255350568308558
$ws.Range("D25").Value = 6.394170072921409
$ws.Range("E25").Value = 9.858471375028445
$ws.Range("F25").Value = 36.18266674523085
$ws.Range("G25").Value = 3.663534593594501
$ws.Range("I25").Value = 28.84510212634324
$ws.Range("K25").Value = 16.01079906252768
$ws.Range("M25").Value = 16.85086998638802

